# DDT(get Data from Excel)-HR
# Rebuilds the CompanyHoliday sheet with a two-column-pair layout and adds
# three new worksheets (CommentTemplates, CriteriaLibrary, EmployeeSuccession)
# with their own seed data, mirroring the authored workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: CompanyHoliday
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Give columns B and D a literal "text" number format *before* writing the
# date-like strings into them, so Excel keeps them as text instead of
# re-interpreting "1/22/2025" etc. as date serials.
$ws1.Columns("B").NumberFormat = "@"
$ws1.Columns("D").NumberFormat = "@"

# Write cells in the same order the original author did, so shared-string
# ids land the same way they did in the authored file.
$ws1.Range("A2").Value = "National Holiday"
$ws1.Range("B2").Value = "1/22/2025"
$ws1.Range("C1").Value = "Holiday Name2"
$ws1.Range("B1").Value = "Date1"
$ws1.Range("A1").Value = "Holiday Name1"
$ws1.Range("D1").Value = "Date2"
$ws1.Range("C2").Value = "New Year"
$ws1.Range("D2").Value = "1/1/2025"

# Header row fill (reuse A1's existing themed fill so no duplicate style is
# created in the process).
$ws1.Range("A1").Copy()
$ws1.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column widths (character units as exposed by the COM layer).
$ws1.Columns("B").ColumnWidth = 15.5
$ws1.Columns("C").ColumnWidth = 19.666666666666668
$ws1.Columns("D").ColumnWidth = 23.333333333333336

# Portrait page orientation.
$ws1.PageSetup.Orientation = 1

# Selection / view state.
$null = $ws1.Range("E7").Select()

# ---------------------------------------------------------------------------
# Sheet 2: CommentTemplates
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "CommentTemplates"

$ws2.Range("A1").Value = "Comment Name"
$ws2.Range("A2").Value = "Performance"
$ws2.Range("B1").Value = "Comments"
$ws2.Range("B2").Value = "Test Comments(Added by Automation)"

$ws1.Range("A1").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Columns("A").ColumnWidth = 17.666666666666668
$ws2.Columns("B").ColumnWidth = 43.666666666666664

$null = $ws2.Range("G19").Select()

# ---------------------------------------------------------------------------
# Sheet 3: CriteriaLibrary (empty placeholder sheet)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "CriteriaLibrary"

$null = $ws3.Range("J15").Select()

# ---------------------------------------------------------------------------
# Sheet 4: EmployeeSuccession
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "EmployeeSuccession"

$ws4.Range("A1").Value = "Employee"
$ws4.Range("B1").Value = "Succession"
$ws4.Range("C1").Value = "Rediness"
$ws4.Range("D1").Value = "Note"
$ws4.Range("A2").Value = "Ali"
$ws4.Range("B2").Value = "amy"
$ws4.Range("C2").Value = "re"
$ws4.Range("D2").Value = "Testing add note succession(Added by Automation)"

$ws1.Range("A1").Copy()
$ws4.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws4.Rows("1").RowHeight = 13.5

$ws4.Columns("A").ColumnWidth = 17.666666666666668
$ws4.Columns("B").ColumnWidth = 25.666666666666668
$ws4.Columns("C").ColumnWidth = 14.333333333333332
$ws4.Columns("D").ColumnWidth = 58.666666666666664

$null = $ws4.Range("D11").Select()

# EmployeeSuccession is the tab that ends up active/selected, matching the
# authored workbook (activeTab points at the last-added sheet).
$ws4.Activate()
